$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.730.18'
$ws.Range('E2').Value = '  +1.12%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.078.59'
$ws.Range('E3').Value = '  +0.79%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.63'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '57.04'
$ws.Range('E8').Value = '  +0.51%  '
$ws.Range('E9').Value = '  +1.24%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0785'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('E11').Value = '  +2.63%  '
$ws.Range('E12').Value = '  +0.19%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.48'
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '20.95'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.760'
$ws.Range('E15').Value = '  -1.65%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.25'
$ws.Range('E16').Value = '  +2.01%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.076.50'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.666.71'
$ws.Range('E18').Value = '  +1.07%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.19'
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '70.75'
$ws.Range('E20').Value = '  +1.84%  '
$ws.Range('E21').Value = '  +0.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '227.79'
$ws.Range('E22').Value = '  +0.74%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('E24').Value = '  -0.86%  '
$ws.Range('E25').Value = '  -0.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.89'
$ws.Range('E26').Value = '  +2.20%  '
$ws.Range('E27').Value = '  +10.92%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.91'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '19.37'
$ws.Range('E30').Value = '  +1.89%  '
$ws.Range('E31').Value = '  +0.48%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.64'
$ws.Range('E32').Value = '  +3.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0625'
$ws.Range('E33').Value = '  +1.32%  '
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.50'
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('E36').Value = '  +2.68%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.37'
$ws.Range('E37').Value = '  +5.19%  '
$ws.Range('E38').Value = '  +0.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.42'
$ws.Range('E39').Value = '  -2.36%  '
$ws.Range('E40').Value = '  +6.01%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.97'
$ws.Range('E41').Value = '  +3.14%  '
$ws.Range('B42').Value = 'HuobiToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.94'
$ws.Range('E42').Value = '  -0.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.39'
$ws.Range('E43').Value = '  +2.76%  '
$ws.Range('E44').Value = '  +1.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.456.61'
$ws.Range('E45').Value = '  -0.96%  '
$ws.Range('E46').Value = '  -0.83%  '
$ws.Range('E47').Value = '  +2.97%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '15.62'
$ws.Range('E48').Value = '  +3.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.42'
$ws.Range('E49').Value = '  +3.52%  '
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '47.41'
$ws.Range('E51').Value = '  +7.48%  '
